$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 1416.25
$ws.Range("J42").Value = 230
$ws.Range("L42").Value = 690
$ws.Range("N42").Value = -1150
$ws.Range("H70").Value = 3083.25
$ws.Range("I70").Value = 3181.7273
$ws.Range("K70").Value = 9545.1819
$ws.Range("M70").Value = -9275.1819
$ws.Range("H73").Value = 3083.25
$ws.Range("I73").Value = 3181.7273
$ws.Range("K73").Value = 9545.1819
$ws.Range("M73").Value = -8609.1819
$ws.Range("H103").Value = 1166.6666
$ws.Range("J103").Value = 1500
$ws.Range("L103").Value = 4500
$ws.Range("N103").Value = -5672
$ws.Range("H113").Value = 5863
$ws.Range("I113").Value = 6150.6665
$ws.Range("K113").Value = 6150.6665
$ws.Range("M113").Value = -2896.6665
$ws.Range("H137").Value = 4065
$ws.Range("J137").Value = 4463.3335
$ws.Range("L137").Value = 13390.0005
$ws.Range("N137").Value = -18490.0005
$ws.Range("H138").Value = 2547.6924
$ws.Range("J138").Value = 3421
$ws.Range("L138").Value = 10263
$ws.Range("N138").Value = -20543

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6695
$ws.Range("I45").Value = 2750
$ws.Range("J45").Value = 9325
$ws.Range("K45").Value = 2750
$ws.Range("L45").Value = 9325
$ws.Range("M45").Value = -2373
$ws.Range("N45").Value = -10079
$ws.Range("H74").Value = 1528.8125
$ws.Range("I74").Value = 1431
$ws.Range("K74").Value = 1431
$ws.Range("M74").Value = -557
$ws.Range("H77").Value = 1528.8125
$ws.Range("I77").Value = 1431
$ws.Range("K77").Value = 7155
$ws.Range("M77").Value = -2787
$ws.Range("H92").Value = 199999
$ws.Range("J92").Value = 199999
$ws.Range("L92").Value = 199999
$ws.Range("N92").Value = -204991
$ws.Range("H122").Value = 8999.5
$ws.Range("I122").Value = 15000
$ws.Range("K122").Value = 45000
$ws.Range("M122").Value = -42550

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H28").Value = 60069
$ws.Range("I28").Value = 60069
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 60069
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -59775
$ws.Range("N28").ClearContents()
$ws.Range("H92").Value = 100000
$ws.Range("J92").Value = 100000
$ws.Range("L92").Value = 100000
$ws.Range("N92").Value = -104992
$ws.Range("H107").Value = 52716.375
$ws.Range("I107").Value = 81339
$ws.Range("K107").Value = 81339
$ws.Range("M107").Value = -79419
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 15680
$ws.Range("J28").Value = 15680
$ws.Range("L28").Value = 15680
$ws.Range("N28").Value = -16170
$ws.Range("H31").Value = 917.5
$ws.Range("I31").Value = 935.8570999999999
$ws.Range("K31").Value = 935.8570999999999
$ws.Range("M31").Value = -640.8570999999999
$ws.Range("H34").Value = 917.5
$ws.Range("I34").Value = 935.8570999999999
$ws.Range("K34").Value = 935.8570999999999
$ws.Range("M34").Value = -733.8570999999999
$ws.Range("H43").Value = 36071.285
$ws.Range("J43").Value = 36071.285
$ws.Range("L43").Value = 36071.285
$ws.Range("N43").Value = -36439.285
$ws.Range("H101").Value = 36071.285
$ws.Range("J101").Value = 36071.285
$ws.Range("L101").Value = 36071.285
$ws.Range("N101").Value = -42561.285
$ws.Range("H122").Value = 650
$ws.Range("I122").Value = 650
$ws.Range("K122").Value = 1950
$ws.Range("M122").Value = 500
$ws.Range("H132").Value = 2129.0625
$ws.Range("J132").Value = 1865
$ws.Range("L132").Value = 5595
$ws.Range("N132").Value = -10655
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 1066.6666
$ws.Range("I64").Value = 1066.6666
$ws.Range("K64").Value = 3199.9998
$ws.Range("M64").Value = -2929.9998
$ws.Range("H67").Value = 1066.6666
$ws.Range("I67").Value = 1066.6666
$ws.Range("K67").Value = 3199.9998
$ws.Range("M67").Value = -2263.9998
$ws.Range("H75").Value = 2115.125
$ws.Range("I75").Value = 150
$ws.Range("K75").Value = 450
$ws.Range("M75").Value = 548
$ws.Range("H78").Value = 2115.125
$ws.Range("I78").Value = 150
$ws.Range("K78").Value = 1350
$ws.Range("M78").Value = 3642
$ws.Range("H92").Value = 658
$ws.Range("I92").Value = 658
$ws.Range("K92").Value = 1974
$ws.Range("M92").Value = -726
$ws.Range("H94").Value = 2697.5
$ws.Range("I94").Value = 395
$ws.Range("J94").Value = 5000
$ws.Range("K94").Value = 1185
$ws.Range("L94").Value = 15000
$ws.Range("M94").Value = -509
$ws.Range("N94").Value = -16352

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 7255.6
$ws.Range("I41").Value = 8993
$ws.Range("J41").Value = 4649.5
$ws.Range("K41").Value = 8993
$ws.Range("L41").Value = 4649.5
$ws.Range("M41").Value = -8638
$ws.Range("N41").Value = -5359.5
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 5500.3335
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H126").Value = 3399.5
$ws.Range("I126").Value = 3399.5
$ws.Range("K126").Value = 10198.5
$ws.Range("M126").Value = -7728.5
$ws.Range("H132").Value = 1024
$ws.Range("I132").Value = 1024
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3072
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -542
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 910.2
$ws.Range("I22").Value = 725.5
$ws.Range("J22").Value = 1033.3334
$ws.Range("K22").Value = 725.5
$ws.Range("L22").Value = 1033.3334
$ws.Range("M22").Value = -430.5
$ws.Range("N22").Value = -1623.3334
$ws.Range("H27").Value = 910.2
$ws.Range("I27").Value = 725.5
$ws.Range("J27").Value = 1033.3334
$ws.Range("K27").Value = 725.5
$ws.Range("L27").Value = 1033.3334
$ws.Range("M27").Value = -618.5
$ws.Range("N27").Value = -1247.3334
$ws.Range("H40").Value = 718030.5600000001
$ws.Range("I40").Value = 4241.8
$ws.Range("J40").Value = 2502502.5
$ws.Range("K40").Value = 4241.8
$ws.Range("L40").Value = 2502502.5
$ws.Range("M40").Value = -4105.8
$ws.Range("N40").Value = -2502774.5
$ws.Range("H55").Value = 590.5
$ws.Range("I55").Value = 590.5
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 590.5
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -417.5
$ws.Range("N55").ClearContents()
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H132").Value = 1399
$ws.Range("I132").Value = 1399
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4197
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1667
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 798312.7
$ws.Range("I136").Value = 835167.3
$ws.Range("J136").Value = 761458
$ws.Range("K136").Value = 2505501.9
$ws.Range("L136").Value = 2284374
$ws.Range("M136").Value = -2502951.9
$ws.Range("N136").Value = -2289474

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 949.5
$ws.Range("I132").Value = 949.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2848.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -318.5
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 1384.2778
$ws.Range("I136").Value = 1355.75
$ws.Range("J136").Value = 1612.5
$ws.Range("K136").Value = 4067.25
$ws.Range("L136").Value = 4837.5
$ws.Range("M136").Value = -1517.25
$ws.Range("N136").Value = -9937.5
